$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44504
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 200
$ws.Range("K2").Value = 800
$ws.Range("L2").Value = 800
$ws.Range("M2").Value = 800
$ws.Range("N2").Value = "$/kilo (volumen en unidades)"
$ws.Range("O2").Value = "Perú"
$ws.Range("P2").Value = 800

# Row 3
$ws.Range("D3").Value = 44477
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 80
$ws.Range("K3").Value = 800
$ws.Range("L3").Value = 800
$ws.Range("M3").Value = 800
$ws.Range("N3").Value = "$/kilo (volumen en unidades)"
$ws.Range("O3").Value = "Perú"
$ws.Range("P3").Value = 800

# Row 4
$ws.Range("D4").Value = 44223
$ws.Range("H4").Value = "Americana O Klondike"
$ws.Range("I4").Value = "Extra"
$ws.Range("J4").Value = 340
$ws.Range("K4").Value = 2500
$ws.Range("L4").Value = 2500
$ws.Range("M4").Value = 2500
$ws.Range("N4").Value = "$/unidad"
$ws.Range("O4").Value = "Región de O'Higgins"
$ws.Range("P4").Value = 2500

# Row 5
$ws.Range("D5").Value = 44223
$ws.Range("H5").Value = "Americana O Klondike"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 400
$ws.Range("K5").Value = 2000
$ws.Range("L5").Value = 2000
$ws.Range("M5").Value = 2000
$ws.Range("N5").Value = "$/unidad"
$ws.Range("O5").Value = "Región de O'Higgins"
$ws.Range("P5").Value = 2000

# Row 6
$ws.Range("D6").Value = 44223
$ws.Range("H6").Value = "Americana O Klondike"
$ws.Range("I6").Value = "Segunda"
$ws.Range("J6").Value = 300
$ws.Range("K6").Value = 1500
$ws.Range("L6").Value = 1500
$ws.Range("M6").Value = 1500
$ws.Range("N6").Value = "$/unidad"
$ws.Range("O6").Value = "Región de O'Higgins"
$ws.Range("P6").Value = 1500

# Row 7
$ws.Range("D7").Value = 44223
$ws.Range("H7").Value = "Americana O Klondike"
$ws.Range("I7").Value = "Tercera"
$ws.Range("J7").Value = 160
$ws.Range("K7").Value = 1000
$ws.Range("L7").Value = 1000
$ws.Range("M7").Value = 1000
$ws.Range("N7").Value = "$/unidad"
$ws.Range("O7").Value = "Región de O'Higgins"
$ws.Range("P7").Value = 1000

# Row 8
$ws.Range("D8").Value = 44510
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 250
$ws.Range("K8").Value = 800
$ws.Range("L8").Value = 800
$ws.Range("M8").Value = 800
$ws.Range("N8").Value = "$/kilo (volumen en unidades)"
$ws.Range("O8").Value = "Perú"
$ws.Range("P8").Value = 800

# Row 9
$ws.Range("D9").Value = 44491
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 150
$ws.Range("K9").Value = 800
$ws.Range("L9").Value = 800
$ws.Range("M9").Value = 800
$ws.Range("N9").Value = "$/kilo (volumen en unidades)"
$ws.Range("O9").Value = "Perú"
$ws.Range("P9").Value = 800

# Row 10
$ws.Range("D10").Value = 44483
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 120
$ws.Range("K10").Value = 800
$ws.Range("L10").Value = 800
$ws.Range("M10").Value = 800
$ws.Range("N10").Value = "$/kilo (volumen en unidades)"
$ws.Range("O10").Value = "Perú"
$ws.Range("P10").Value = 800

# Row 11
$ws.Range("D11").Value = 44312
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 180
$ws.Range("K11").Value = 2500
$ws.Range("L11").Value = 2500
$ws.Range("M11").Value = 2500
$ws.Range("N11").Value = "$/unidad"
$ws.Range("O11").Value = "Perú"
$ws.Range("P11").Value = 2500

# Row 12
$ws.Range("D12").Value = 44167
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 400
$ws.Range("K12").Value = 5000
$ws.Range("L12").Value = 5000
$ws.Range("M12").Value = 5000
$ws.Range("N12").Value = "$/unidad"
$ws.Range("O12").Value = "Región de O'Higgins"
$ws.Range("P12").Value = 5000

# Row 13
$ws.Range("D13").Value = 44167
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Segunda"
$ws.Range("J13").Value = 560
$ws.Range("K13").Value = 3000
$ws.Range("L13").Value = 3000
$ws.Range("M13").Value = 3000
$ws.Range("N13").Value = "$/unidad"
$ws.Range("O13").Value = "Región de O'Higgins"
$ws.Range("P13").Value = 3000

# Row 14
$ws.Range("D14").Value = 44167
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Tercera"
$ws.Range("J14").Value = 450
$ws.Range("K14").Value = 2000
$ws.Range("L14").Value = 2000
$ws.Range("M14").Value = 2000
$ws.Range("N14").Value = "$/unidad"
$ws.Range("O14").Value = "Región de O'Higgins"
$ws.Range("P14").Value = 2000

# Row 15
$ws.Range("D15").Value = 44305
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 100
$ws.Range("K15").Value = 2500
$ws.Range("L15").Value = 2500
$ws.Range("M15").Value = 2500
$ws.Range("N15").Value = "$/unidad"
$ws.Range("O15").Value = "Perú"
$ws.Range("P15").Value = 2500

# Row 16
$ws.Range("D16").Value = 44497
$ws.Range("H16").Value = "Sin especificar"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 250
$ws.Range("K16").Value = 800
$ws.Range("L16").Value = 800
$ws.Range("M16").Value = 800
$ws.Range("N16").Value = "$/kilo (volumen en unidades)"
$ws.Range("O16").Value = "Perú"
$ws.Range("P16").Value = 800

# Row 17
$ws.Range("D17").Value = 44488
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 150
$ws.Range("K17").Value = 800
$ws.Range("L17").Value = 800
$ws.Range("M17").Value = 800
$ws.Range("N17").Value = "$/kilo (volumen en unidades)"
$ws.Range("O17").Value = "Perú"
$ws.Range("P17").Value = 800

# Row 18
$ws.Range("D18").Value = 44495
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 200
$ws.Range("K18").Value = 800
$ws.Range("L18").Value = 800
$ws.Range("M18").Value = 800
$ws.Range("N18").Value = "$/kilo (volumen en unidades)"
$ws.Range("O18").Value = "Perú"
$ws.Range("P18").Value = 800

# Row 19
$ws.Range("D19").Value = 44217
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Extra"
$ws.Range("J19").Value = 400
$ws.Range("K19").Value = 2500
$ws.Range("L19").Value = 2500
$ws.Range("M19").Value = 2500
$ws.Range("N19").Value = "$/unidad"
$ws.Range("O19").Value = "Región de O'Higgins"
$ws.Range("P19").Value = 2500

# Row 20
$ws.Range("D20").Value = 44217
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 280
$ws.Range("K20").Value = 2000
$ws.Range("L20").Value = 2000
$ws.Range("M20").Value = 2000
$ws.Range("N20").Value = "$/unidad"
$ws.Range("O20").Value = "Región de O'Higgins"
$ws.Range("P20").Value = 2000

# Row 21
$ws.Range("D21").Value = 44194
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Extra"
$ws.Range("J21").Value = 120
$ws.Range("K21").Value = 3500
$ws.Range("L21").Value = 3500
$ws.Range("M21").Value = 3500
$ws.Range("N21").Value = "$/unidad"
$ws.Range("O21").Value = "Región de O'Higgins"
$ws.Range("P21").Value = 3500

# Row 22
$ws.Range("D22").Value = 44194
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 3000
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = 3000
$ws.Range("N22").Value = "$/unidad"
$ws.Range("O22").Value = "Región de O'Higgins"
$ws.Range("P22").Value = 3000
